# This file (trials 2..41, i.e. subject_id 39's 40-trial block) is being
# regenerated for a different subject: same 40 stimuli, new trial_total
# numbering (column F, decremented by 162 so the running trial counter is
# continuous with the subject's other blocks) and the per-row stimulus
# payload (H,I,K,L,M,N,O,P,Q,R,S,T,U,V) reshuffled onto a new row order
# ("make only 20 different versions and duplicate many times for 1000
# subjects"). Columns A-E (subject_id/task/block_total/block_scene/
# trial_block), G (target_cat) and J (cond_mem) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 41
$fDelta   = -162

# row(new) -> row(old) : which old row's stimulus payload lands on this row
$rowMap = @{
  2=33; 3=19; 4=28; 5=12; 6=13; 7=17; 8=18; 9=11; 10=3; 11=14;
  12=34; 13=27; 14=30; 15=32; 16=29; 17=2; 18=41; 19=26; 20=23; 21=40;
  22=35; 23=16; 24=37; 25=4; 26=25; 27=21; 28=6; 29=24; 30=20; 31=5;
  32=36; 33=8; 34=10; 35=7; 36=9; 37=39; 38=22; 39=31; 40=38; 41=15
}

# Payload columns copied row-to-row as a block (by column index).
$payloadCols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# Snapshot every old value up front so writes for one row never clobber
# data another row still needs to read (the map above is a single 40-cycle).
# NOTE: use .Value2 (not .Value) to read - .Value's getter is unreliable
# for round-tripping through a variable in this host.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $payloadCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # trial_total just shifts down by a constant.
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 6).Value2 + $fDelta

    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $payloadCols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
